$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Thursday 2021-09-09
$ws.Range("A9").Value = 44448
$ws.Range("B9").Value = "8 hours"
$ws.Range("C9").Value = "Task:"

# Row 10 - Sunday 2021-09-12
$ws.Range("A10").Value = 44451
$ws.Range("A10").NumberFormat = "d-mmm"
$ws.Range("B10").Value = "7 hours"
$ws.Range("C10").Value = "Task: Finished up the task and tested all corner cases, also added multi language"

# Row 11 - Sunday 2021-09-12
$ws.Range("A11").Value = 44451
$ws.Range("A11").NumberFormat = "d-mmm"
$ws.Range("C11").Value = "Self Learning: Started watching the pluralsight asp mvc course"
$ws.Range("B11").Value = "1 hour"

# Row 12 - Monday 2021-09-13
$ws.Range("A12").Value = 44452
$ws.Range("A12").NumberFormat = "d-mmm"
$ws.Range("B12").Value = "8 hours"
$ws.Range("C12").Value = "Task: Addressing the comments I got on the task"

# Row 13 - Monday 2021-09-13
$ws.Range("A13").Value = 44452
$ws.Range("A13").NumberFormat = "d-mmm"
$ws.Range("B13").Value = "1 hour"
$ws.Range("C13").Value = "Self Learning: Continuing with the asp.net mvc course on pluralsight"

# Update the active selection to A14, matching the end-of-day commit's cursor position
$ws.Range("A14").Select()
